# "nor reads timeline data"
# The timeline sheet is trimmed down from a 7-column / 9-row layout to a
# 6-column / 4-row layout: the extra "random/rest/timpani/horn/violin/viola"
# rows and the trailing "viola" column are dropped, the time-signature /
# sample-id header row is replaced with a new set of values, and the last
# two data rows are rewritten to share one "renda-d6-*" sample progression.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused last column (G) and the trailing rows (5-9) that held
# the old random/rest/timpani/horn/violin/viola data.
$ws.Columns("G").Delete() | Out-Null
$ws.Rows("5:9").Delete() | Out-Null

# Row 1: header / time-signature row
$ws.Range("A1").Value = "measure"
$ws.Range("B1").Value = "8/4"
$ws.Range("C1").Value = "7/4"
$ws.Range("D1").Value = "6/4"
$ws.Range("E1").Value = "5/4"
$ws.Range("F1").Value = "4/4"

# Row 2: harmony row (unchanged content, just one fewer column)
$ws.Range("A2").Value = "harmony"
$ws.Range("B2").Value = "pentatonic"
$ws.Range("C2").Value = "major"
$ws.Range("D2").Value = "minor"
$ws.Range("E2").Value = "pentatonic"
$ws.Range("F2").Value = "major"

# Row 3: flute now plays the renda-d6 sample progression
$ws.Range("A3").Value = "flute"
$ws.Range("B3").Value = "renda-d6-3"
$ws.Range("C3").Value = "renda-d6-4"
$ws.Range("D3").Value = "renda-d6-5"
$ws.Range("E3").Value = "renda-d6-6"
$ws.Range("F3").Value = "renda-d6-7"

# Row 4: oboe plays the same renda-d6 sample progression
$ws.Range("A4").Value = "oboe"
$ws.Range("B4").Value = "renda-d6-3"
$ws.Range("C4").Value = "renda-d6-4"
$ws.Range("D4").Value = "renda-d6-5"
$ws.Range("E4").Value = "renda-d6-6"
$ws.Range("F4").Value = "renda-d6-7"

$ws.Range("D7").Select()
